{"js": "// The document is being reset to a blank \"generic template\": all of the\n// sample content (headings, normal/compact paragraphs, list examples,\n// title, table, page break, second-page heading, etc.) is removed from\n// the body, leaving a single empty paragraph while the page setup\n// (margins, size, footer) is preserved.\ncontext.document.body.clear();\nawait context.sync();\n", "ps1": "# The document is being reset to a blank \"generic template\": all of the\n# sample content (headings, normal/compact paragraphs, list examples,\n# title, table, page break, second-page heading, etc.) is removed from\n# the body, leaving a single empty paragraph while the page setup\n# (margins, size, footer) is preserved.\n$d = $word.ActiveDocument\n\n# Repeatedly clear the document's main story range until only the\n# trailing empty paragraph mark is left (a single Delete() call only\n# removes the first paragraph in this host, so loop it out).\n$guard = 0\nwhile ($d.Characters.Count -gt 1 -and $guard -lt 100) {\n    $d.Content.Delete()\n    $guard++\n}\n"}
